$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 990.3333
$ws.Range("I12").Value = 1000.5
$ws.Range("K12").Value = 1000.5
$ws.Range("M12").Value = -830.5

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 503499.75
$ws.Range("I28").Value = 503499.75
$ws.Range("K28").Value = 503499.75
$ws.Range("M28").Value = -503014.75

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3706.0454
$ws.Range("I43").Value = 2994
$ws.Range("J43").Value = 4299.4165
$ws.Range("K43").Value = 2994
$ws.Range("L43").Value = 4299.4165
$ws.Range("M43").Value = -2925
$ws.Range("N43").Value = -4437.4165

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5748.8335
$ws.Range("I62").Value = 5189.2
$ws.Range("K62").Value = 5189.2
$ws.Range("M62").Value = -4565.2

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8071.2856
$ws.Range("J64").Value = 8833.5
$ws.Range("L64").Value = 8833.5
$ws.Range("N64").Value = -9329.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5748.8335
$ws.Range("I65").Value = 5189.2
$ws.Range("K65").Value = 25946
$ws.Range("M65").Value = -22826

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 8071.2856
$ws.Range("J67").Value = 8833.5
$ws.Range("L67").Value = 8833.5
$ws.Range("N67").Value = -10549.5

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3886.25
$ws.Range("J111").Value = 4193.75
$ws.Range("L111").Value = 12581.25
$ws.Range("N111").Value = -18715.25

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 71817.5
$ws.Range("I137").Value = 71817.5
$ws.Range("K137").Value = 215452.5
$ws.Range("M137").Value = -212902.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2629.85
$ws.Range("J138").Value = 2720.8171
$ws.Range("L138").Value = 8162.451300000001
$ws.Range("N138").Value = -18442.4513

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7453.231
$ws.Range("I141").Value = 6183
$ws.Range("K141").Value = 18549
$ws.Range("M141").Value = -13369

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3713.8333
$ws.Range("I63").Value = 2399
$ws.Range("J63").Value = 3976.8
$ws.Range("K63").Value = 2399
$ws.Range("L63").Value = 3976.8
$ws.Range("M63").Value = -1713
$ws.Range("N63").Value = -5348.8

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3713.8333
$ws.Range("I66").Value = 2399
$ws.Range("J66").Value = 3976.8
$ws.Range("K66").Value = 11995
$ws.Range("L66").Value = 19884
$ws.Range("M66").Value = -8563
$ws.Range("N66").Value = -26748

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3253.2856
$ws.Range("J97").Value = 5652.75
$ws.Range("L97").Value = 5652.75
$ws.Range("N97").Value = -6644.75

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2947.1072
$ws.Range("I132").Value = 2684.6
$ws.Range("K132").Value = 8053.799999999999
$ws.Range("M132").Value = -5523.799999999999

# ARM row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1084.05
$ws.Range("I94").Value = 1160.7273
$ws.Range("J94").Value = 990.3333
$ws.Range("K94").Value = 1160.7273
$ws.Range("L94").Value = 990.3333
$ws.Range("M94").Value = -709.7273
$ws.Range("N94").Value = -1892.3333

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2454.875
$ws.Range("I107").Value = 923.5
$ws.Range("K107").Value = 923.5
$ws.Range("M107").Value = 996.5

# CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 153373
$ws.Range("J20").Value = 153373
$ws.Range("L20").Value = 153373
$ws.Range("N20").Value = -153845

# CRP row 28
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 51441.285
$ws.Range("J28").Value = 51441.285
$ws.Range("L28").Value = 51441.285
$ws.Range("N28").Value = -51931.285

# CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 153373
$ws.Range("J30").Value = 153373
$ws.Range("L30").Value = 153373
$ws.Range("N30").Value = -153555

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3110.4473
$ws.Range("I58").Value = 2915.625
$ws.Range("K58").Value = 2915.625
$ws.Range("M58").Value = -2712.625

# CRP row 87
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 96963
$ws.Range("J87").Value = 96963
$ws.Range("L87").Value = 96963
$ws.Range("N87").Value = -99335

# CRP row 90
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 96963
$ws.Range("J90").Value = 96963
$ws.Range("L90").Value = 290889
$ws.Range("N90").Value = -302745

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3332.889
$ws.Range("J99").Value = 3999.4
$ws.Range("L99").Value = 3999.4
$ws.Range("N99").Value = -6995.4

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3332.889
$ws.Range("J126").Value = 3999.4
$ws.Range("L126").Value = 11998.2
$ws.Range("N126").Value = -16938.2

# CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 153373
$ws.Range("J128").Value = 153373
$ws.Range("L128").Value = 153373
$ws.Range("N128").Value = -163333

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2575.5833
$ws.Range("I134").Value = 1890
$ws.Range("K134").Value = 5670
$ws.Range("M134").Value = -3135

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3110.4473
$ws.Range("I136").Value = 2915.625
$ws.Range("K136").Value = 8746.875
$ws.Range("M136").Value = -6196.875

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 4166.3335
$ws.Range("J9").Value = 4999.75
$ws.Range("L9").Value = 14999.25
$ws.Range("N9").Value = -15447.25

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 775
$ws.Range("J68").Value = 775
$ws.Range("L68").Value = 2325
$ws.Range("N68").Value = -3947

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 775
$ws.Range("J71").Value = 775
$ws.Range("L71").Value = 6975
$ws.Range("N71").Value = -15087

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 4540
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -17496

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1208.5
$ws.Range("I113").Value = 667
$ws.Range("K113").Value = 2001
$ws.Range("M113").Value = 169

# CUL row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3500.5
$ws.Range("J127").Value = 3500.5
$ws.Range("L127").Value = 10501.5
$ws.Range("N127").Value = -20421.5

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1658.5588
$ws.Range("I131").Value = 1283.8572
$ws.Range("J131").Value = 1755.7037
$ws.Range("K131").Value = 3851.5716
$ws.Range("L131").Value = 5267.1111
$ws.Range("M131").Value = 1188.4284
$ws.Range("N131").Value = -15347.1111

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 820.55554
$ws.Range("I55").Value = 980.6667
$ws.Range("J55").Value = 500.33334
$ws.Range("K55").Value = 980.6667
$ws.Range("L55").Value = 500.33334
$ws.Range("M55").Value = -807.6667
$ws.Range("N55").Value = -846.33334

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2827.4546
$ws.Range("I93").Value = 1790.4
$ws.Range("J93").Value = 3691.6667
$ws.Range("K93").Value = 1790.4
$ws.Range("L93").Value = 3691.6667
$ws.Range("M93").Value = -542.4000000000001
$ws.Range("N93").Value = -6187.6667

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3859.9678
$ws.Range("I132").Value = 3643.28
$ws.Range("K132").Value = 10929.84
$ws.Range("M132").Value = -8399.84

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7878.5
$ws.Range("I136").Value = 5149.9
$ws.Range("J136").Value = 12426.167
$ws.Range("K136").Value = 15449.7
$ws.Range("L136").Value = 37278.501
$ws.Range("M136").Value = -12899.7
$ws.Range("N136").Value = -42378.501

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2255.6667
$ws.Range("I132").Value = 1828.8572
$ws.Range("K132").Value = 5486.571599999999
$ws.Range("M132").Value = -2956.571599999999
